$d = $word.ActiveDocument

# Locate the first "Completed tasks: " paragraph (the Cassandra implementation
# section) by searching for the distinctive text and anchoring on its range.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Completed tasks:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Expand to the full paragraph so we can insert a new paragraph right after it.
$para = $rng.Paragraphs(1)
$paraRange = $para.Range

$paraRange.InsertParagraphAfter()

# The freshly inserted paragraph mark creates a new paragraph; grab it via
# the paragraph that now follows the "Completed tasks:" paragraph.
$newPara = $para.Next()
$newRange = $newPara.Range
$newRange.Text = "Configuration of Cassandra on compute nodes"

$newRange.Font.Name = "Georgia"
$newRange.Font.Size = 10

$newPara.Range.ListFormat.ApplyNumberDefault()
